$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.3295
$ws.Range("D3").Value = -7.00319999999999
$ws.Range("E8").Value = 16.22930000000001
$ws.Range("E11").Value = 16.5701
$ws.Range("A12").Value = -21.57050000000001
$ws.Range("B14").Value = 6.887299999999999
$ws.Range("E14").Value = 16.8105
$ws.Range("E15").Value = 16.4207
$ws.Range("B26").Value = 3.746000000000005
$ws.Range("D30").Value = -7.801600000000005
$ws.Range("B31").Value = 4.754800000000006
$ws.Range("A32").Value = -21.27330000000001
$ws.Range("B35").Value = 9.017000000000003
$ws.Range("A36").Value = -19.8796
$ws.Range("E36").Value = 15.9937
$ws.Range("B37").Value = 8.826600000000004
$ws.Range("A38").Value = -19.4142
$ws.Range("D44").Value = -7.132399999999999
$ws.Range("B45").Value = 5.9072
$ws.Range("A46").Value = -21.6188
$ws.Range("A54").Value = -21.86049999999999
$ws.Range("A55").Value = -22.51490000000002
$ws.Range("B57").Value = 4.878599999999995
$ws.Range("D58").Value = -8.221499999999997
$ws.Range("E64").Value = 17.4199
$ws.Range("A67").Value = -21.47579999999998
$ws.Range("A69").Value = -21.61579999999997
$ws.Range("A72").Value = -22.12660000000001
$ws.Range("D84").Value = -8.748900000000003
$ws.Range("D89").Value = -5.981799999999997
$ws.Range("E89").Value = 18.62000000000003
$ws.Range("A91").Value = -21.36110000000001
$ws.Range("D91").Value = -5.988399999999997
$ws.Range("D92").Value = -5.967799999999995
$ws.Range("A99").Value = -20.14109999999998
$ws.Range("B100").Value = 5.303199999999997
$ws.Range("B102").Value = 8.269400000000005
$ws.Range("D102").Value = -8.0052
